# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# for specific leve rows across multiple job sheets, per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 5499.5
$ws.Range("I20").Value = 1999
$ws.Range("J20").Value = 9000
$ws.Range("K20").Value = 1999
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -1769
$ws.Range("N20").Value = -9460
# Row 32
$ws.Range("H32").Value = 3224.2856
$ws.Range("I32").Value = 3595
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 3595
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -3269
$ws.Range("N32").Value = -1652
# Row 35
$ws.Range("H35").Value = 5499.5
$ws.Range("I35").Value = 1999
$ws.Range("J35").Value = 9000
$ws.Range("K35").Value = 1999
$ws.Range("L35").Value = 9000
$ws.Range("M35").Value = -1620
$ws.Range("N35").Value = -9758
# Row 125
$ws.Range("H125").Value = 1301.5555
$ws.Range("I125").Value = 1287.7142
$ws.Range("J125").Value = 1350
$ws.Range("K125").Value = 11589.4278
$ws.Range("L125").Value = 12150
$ws.Range("M125").Value = -9129.427799999999
$ws.Range("N125").Value = -17070
# Row 138
$ws.Range("H138").Value = 2443.6453
$ws.Range("I138").Value = 1827.2858
$ws.Range("J138").Value = 3242.6296
$ws.Range("K138").Value = 5481.857400000001
$ws.Range("L138").Value = 9727.888800000001
$ws.Range("M138").Value = -341.8574000000008
$ws.Range("N138").Value = -20007.8888

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 4897.9375
$ws.Range("I74").Value = 1126
$ws.Range("J74").Value = 6373.913
$ws.Range("K74").Value = 1126
$ws.Range("L74").Value = 6373.913
$ws.Range("M74").Value = -252
$ws.Range("N74").Value = -8121.913
# Row 77
$ws.Range("H77").Value = 4897.9375
$ws.Range("I77").Value = 1126
$ws.Range("J77").Value = 6373.913
$ws.Range("K77").Value = 5630
$ws.Range("L77").Value = 31869.565
$ws.Range("M77").Value = -1262
$ws.Range("N77").Value = -40605.565
# Row 110
$ws.Range("H110").Value = 1392.9333
$ws.Range("I110").Value = 1191.4
$ws.Range("K110").Value = 1191.4
$ws.Range("M110").Value = 853.5999999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 562.05554
$ws.Range("I5").Value = 271.8
$ws.Range("J5").Value = 673.6923
$ws.Range("K5").Value = 271.8
$ws.Range("L5").Value = 673.6923
$ws.Range("M5").Value = -159.8
$ws.Range("N5").Value = -897.6923
# Row 10
$ws.Range("H10").Value = 1162.6
$ws.Range("I10").Value = 1543.7142
$ws.Range("J10").Value = 273.33334
$ws.Range("K10").Value = 1543.7142
$ws.Range("L10").Value = 273.33334
$ws.Range("M10").Value = -1404.7142
$ws.Range("N10").Value = -551.33334
# Row 25
$ws.Range("H25").Value = 1083
$ws.Range("I25").Value = 1083
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1083
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -909
$ws.Range("N25").ClearContents()
# Row 31
$ws.Range("H31").Value = 1320.3704
$ws.Range("J31").Value = 1376.5
$ws.Range("L31").Value = 1376.5
$ws.Range("N31").Value = -1966.5
# Row 34
$ws.Range("H34").Value = 1320.3704
$ws.Range("J34").Value = 1376.5
$ws.Range("L34").Value = 1376.5
$ws.Range("N34").Value = -1780.5
# Row 39
$ws.Range("H39").Value = 10051
$ws.Range("I39").Value = 10051
$ws.Range("K39").Value = 10051
$ws.Range("M39").Value = -9660
# Row 49
$ws.Range("H49").Value = 10051
$ws.Range("I49").Value = 10051
$ws.Range("K49").Value = 10051
$ws.Range("M49").Value = -9869
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 30304350
$ws.Range("J122").Value = 3816.3333
$ws.Range("L122").Value = 34346.9997
$ws.Range("N122").Value = -39246.9997
# Row 131
$ws.Range("H131").Value = 1494808.1
$ws.Range("I131").Value = 4488.3335
$ws.Range("J131").Value = 1819968.8
$ws.Range("K131").Value = 13465.0005
$ws.Range("L131").Value = 5459906.4
$ws.Range("M131").Value = -8425.000499999998
$ws.Range("N131").Value = -5469986.4

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 23000
$ws.Range("J39").Value = 23000
$ws.Range("L39").Value = 23000
$ws.Range("N39").Value = -24064
# Row 74
$ws.Range("H74").Value = 59800
$ws.Range("J74").Value = 59800
$ws.Range("L74").Value = 59800
$ws.Range("N74").Value = -61672
# Row 77
$ws.Range("H77").Value = 59800
$ws.Range("J77").Value = 59800
$ws.Range("L77").Value = 179400
$ws.Range("N77").Value = -188760
# Row 126
$ws.Range("H126").Value = 1702
$ws.Range("I126").Value = 1632
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 4896
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2426
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1417.1154
$ws.Range("I7").Value = 1179.3334
$ws.Range("J7").Value = 1741.3636
$ws.Range("K7").Value = 1179.3334
$ws.Range("L7").Value = 1741.3636
$ws.Range("M7").Value = -1067.3334
$ws.Range("N7").Value = -1965.3636
# Row 31
$ws.Range("H31").Value = 3050
$ws.Range("I31").Value = 200
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = 48
$ws.Range("N31").Value = -4496
# Row 122
$ws.Range("H122").Value = 6844.857
$ws.Range("I122").Value = 7500.0527
$ws.Range("J122").Value = 5461.6665
$ws.Range("K122").Value = 22500.1581
$ws.Range("L122").Value = 16384.9995
$ws.Range("M122").Value = -20050.1581
$ws.Range("N122").Value = -21284.9995
# Row 126
$ws.Range("H126").Value = 1417.1154
$ws.Range("I126").Value = 1179.3334
$ws.Range("J126").Value = 1741.3636
$ws.Range("K126").Value = 3538.0002
$ws.Range("L126").Value = 5224.0908
$ws.Range("M126").Value = -1068.0002
$ws.Range("N126").Value = -10164.0908

$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 4400
$ws.Range("I32").Value = 4400
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4400
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4083
$ws.Range("N32").ClearContents()
# Row 63
$ws.Range("H63").Value = 36665.332
$ws.Range("J63").Value = 36665.332
$ws.Range("L63").Value = 36665.332
$ws.Range("N63").Value = -37913.332
# Row 66
$ws.Range("H66").Value = 36665.332
$ws.Range("J66").Value = 36665.332
$ws.Range("L66").Value = 109995.996
$ws.Range("N66").Value = -116235.996
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 75
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36872
# Row 78
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -114360
# Row 122
$ws.Range("H122").Value = 1500.1177
$ws.Range("I122").Value = 1483.8182
$ws.Range("J122").Value = 1530
$ws.Range("K122").Value = 4451.4546
$ws.Range("M122").Value = -2001.4546
$ws.Range("N122").Value = -9490
